$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '42.857.72'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +1.20%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.301.02'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.41%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '316.42'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.04%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '103.94'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  +0.25%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.601'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -1.88%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '39.31'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -1.76%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0906'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.48%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '8.47'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('E14').Value = '  +4.27%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '15.30'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.05%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.651.46'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +0.13%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.301.23'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.06%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '42.810.11'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.48'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '14.16'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +27.01%  '
$ws.Range('E21').Value = '  -0.55%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '74.00'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +1.19%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '3.56'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.21%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '264.57'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -4.54%  '
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E26').Value = '  +0.60%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '10.88'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('E28').Value = '  -0.02%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '7.06'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +20.14%  '
$ws.Range('E30').Value = '  -1.55%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '37.57'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +4.83%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '166.91'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.97%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.0872'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  -4.54%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('E37').Value = '  -0.37%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0350'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -5.75%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.76'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  -2.94%  '
$ws.Range('E41').Value = '  +4.74%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.230'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +1.47%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '69.14'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('E44').Value = '  +0.27%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '92.17'
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '12.33'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +2.19%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '114.58'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +1.48%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.731.69'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +8.91%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '80.22'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.72%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '8.78'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.37%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '5.14'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +0.74%  '
